$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 301
$ws1.Range("F8").Value = 1603
$ws1.Range("F12").Value = 2633
$ws1.Range("F15").Value = 6915
$ws1.Range("F17").Value = 7062
$ws1.Range("F18").Value = 7062
$ws1.Range("F20").Value = 2670
$ws1.Range("F24").Value = 140
$ws1.Range("F25").Value = 1811
$ws1.Range("F29").Value = 9
$ws1.Range("F31").Value = 33
$ws1.Range("F33").Value = 1082
$ws1.Range("F34").Value = 2517
$ws1.Range("F37").Value = 365
$ws1.Range("F38").Value = 1004
$ws1.Range("F40").Value = 454

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 301
$ws4.Range("F9").Value = 1603
$ws4.Range("F15").Value = 2633
$ws4.Range("F21").Value = 6915
$ws4.Range("F23").Value = 7062
$ws4.Range("F24").Value = 7062
$ws4.Range("F26").Value = 2671
$ws4.Range("F33").Value = 1811
$ws4.Range("F39").Value = 33
$ws4.Range("F42").Value = 2517
$ws4.Range("F46").Value = 365
$ws4.Range("F47").Value = 1004
$ws4.Range("F49").Value = 454
